# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.918.53'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '3.391.57'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'579.96"
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").Value = "'138.14"
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.390.39'
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").Value = "'7.53"
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("D13").Value = '3.972.14'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("D16").Value = '3.397.33'
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").Value = "'25.42"
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = '61.977.19'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").Value = "'14.16"
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D20").Value = "'5.82"
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = "'9.47"
$ws.Range("E21").Value = '  +1.13%  '
$ws.Range("D22").Value = "'392.78"
$ws.Range("E22").Value = '  +1.92%  '
$ws.Range("D23").Value = "'0.566"
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("D24").Value = '3.540.89'
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D25").Value = "'0.0000130"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = "'71.57"
$ws.Range("E27").Value = '  +0.84%  '
$ws.Range("E28").Value = '  -3.14%  '
$ws.Range("E29").Value = '  -2.76%  '
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  +2.72%  '
$ws.Range("D32").Value = "'8.23"
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").Value = '3.424.26'
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("E37").Value = '  -3.80%  '
$ws.Range("E38").Value = '  +2.04%  '
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("D40").Value = "'164.92"
$ws.Range("E40").Value = '  +1.15%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("E42").Value = '  +8.84%  '
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("E44").Value = '  +3.68%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = "'4.44"
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").Value = "'25.04"
$ws.Range("E47").Value = '  +6.39%  '
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").Value = "'6.90"
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").Value = "'23.10"
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").Value = '2.339.26'
$ws.Range("E51").Value = '  +6.28%  '
